$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data rows 2-11: A,B,C,D,E,F,G,H,I,J
$data = @(
    @(0, 0, 565.1923707142779, 45413.90750580061, 45413.907636235,   [double]"0.0001304343865740741", 25, 28, 26, "COMPLETE"),
    @(1, 1, 703.9281400860115, 45413.90764671707, 45413.90765120689, [double]"4.489814814814815e-06", 14, 2,  4,  "COMPLETE"),
    @(2, 2, 744.2480782766445, 45413.90765201986, 45413.90765606792, [double]"4.048055555555555e-06", 22, 10, 13, "COMPLETE"),
    @(3, 3, 633.050001047896,  45413.9076567984,  45413.90766115845, [double]"4.360046296296296e-06", 23, 21, 13, "COMPLETE"),
    @(4, 4, 567.1161762126092, 45413.90766198299, 45413.9076662549,  [double]"4.271909722222223e-06", 25, 27, 15, "COMPLETE"),
    @(5, 5, 706.0405603478761, 45413.90766709617, 45413.90767021354, [double]"3.117372685185185e-06", 19, 1,  30, "COMPLETE"),
    @(6, 6, 631.2513793815023, 45413.90767076387, 45413.90767342736, [double]"2.663483796296296e-06", 26, 25, 21, "COMPLETE"),
    @(7, 7, 643.3065150904984, 45413.90767412398, 45413.90767807258, [double]"3.948599537037037e-06", 10, 11, 21, "COMPLETE"),
    @(8, 8, 600.0418887143586, 45413.90767898317, 45413.90768313344, [double]"4.150266203703703e-06", 8,  12, 17, "COMPLETE"),
    @(9, 9, 744.8583513583845, 45413.90769887075, 45413.90770314227, [double]"4.271516203703704e-06", 15, 23, 19, "COMPLETE")
)

$rowIndex = 2
foreach ($row in $data) {
    $ws.Cells.Item($rowIndex, 1).Value = $row[0]
    $ws.Cells.Item($rowIndex, 2).Value = $row[1]
    $ws.Cells.Item($rowIndex, 3).Value = $row[2]
    $ws.Cells.Item($rowIndex, 4).Value = $row[3]
    $ws.Cells.Item($rowIndex, 5).Value = $row[4]
    $ws.Cells.Item($rowIndex, 6).Value = $row[5]
    $ws.Cells.Item($rowIndex, 7).Value = $row[6]
    $ws.Cells.Item($rowIndex, 8).Value = $row[7]
    $ws.Cells.Item($rowIndex, 9).Value = $row[8]
    $ws.Cells.Item($rowIndex, 10).Value = $row[9]
    $rowIndex++
}

# Copy styles from row 3 template to the newly added rows 4-11
$ws.Range("A3:J3").Copy() | Out-Null
$ws.Range("A4:J11").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

Write-Output "done"
